$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handback identifiers (replacing the old ones throughout the workbook)
# ---------------------------------------------------------------------------
#   a428c49f-18c1-433b-92e0-e4b70ebbbad9  ->  2bf6a623-283a-48ef-9d6e-1f7852ead1d9
#   b3dd4da3-1e78-4023-8ce8-6b067d55c8b1  ->  ffff90e2eb77-fe90-4e88-8774-6b89f0659ca8

$oldId1 = "a428c49f-18c1-433b-92e0-e4b70ebbbad9"
$newId1 = "2bf6a623-283a-48ef-9d6e-1f7852ead1d9"
$oldId2 = "b3dd4da3-1e78-4023-8ce8-6b067d55c8b1"
$newId2 = "ffff90e2eb77-fe90-4e88-8774-6b89f0659ca8"

$newMd1 = "$newId1.md"
$newMd2 = "$newId2.md"

$newZhCnXlf = "$newId1.620d4ce6edac63726341fd8012fba993f62faff2.zh-cn.xlf"
$newDeDeXlf = "$newId1.620d4ce6edac63726341fd8012fba993f62faff2.de-de.xlf"

$latestHoDate  = "2016-08-24 19:12:35"
$zhHandoffDate = "2016-08-24 19:12:30"
$zhHandbackDate= "2016-08-24 19:12:47"
$deHandbackDate= "2016-08-24 19:12:54"

# ===========================================================================
# Overview sheet
# ===========================================================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $newMd1
$ov.Range("A3").Value = $newMd2
$ov.Range("G2").Value = $latestHoDate
$ov.Range("G3").Value = $latestHoDate

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/89a8524106103f1b2796df6de832e7068317d359/e2e/$newMd1", [Type]::Missing, [Type]::Missing, "e2e\$newMd1")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/89a8524106103f1b2796df6de832e7068317d359/e2e/$newMd2", [Type]::Missing, [Type]::Missing, "e2e\$newMd2")

# ===========================================================================
# zh-cn sheet
# ===========================================================================
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $newMd1
$zh.Range("I2").Value = $newMd1
$zh.Range("A3").Value = $newMd2
$zh.Range("I3").Value = $newMd2

$zh.Range("G2").Value = $newZhCnXlf
$zh.Range("H2").Value = $zhHandoffDate
$zh.Range("J2").Value = $newZhCnXlf
$zh.Range("K2").Value = $zhHandbackDate

# Row 3 mirrors row 2's handoff/handback file + datetimes (matches source data)
$zh.Range("G3").Value = $newZhCnXlf
$zh.Range("H3").Value = $zhHandoffDate
$zh.Range("J3").Value = $newZhCnXlf
$zh.Range("K3").Value = $zhHandbackDate

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/89a8524106103f1b2796df6de832e7068317d359/e2e/$newMd1", [Type]::Missing, [Type]::Missing, $newMd1)
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/aa451823352176554acd78a211047dbb3eaa1a45/e2e/$newMd1", [Type]::Missing, [Type]::Missing, $newMd1)
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/89a8524106103f1b2796df6de832e7068317d359/e2e/$newMd2", [Type]::Missing, [Type]::Missing, $newMd2)
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/aa451823352176554acd78a211047dbb3eaa1a45/e2e/$newMd2", [Type]::Missing, [Type]::Missing, $newMd2)

# ===========================================================================
# de-de sheet
# ===========================================================================
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $newMd1
$de.Range("I2").Value = $newMd1
$de.Range("A3").Value = $newMd2
$de.Range("I3").Value = $newMd2

$de.Range("G2").Value = $newDeDeXlf
$de.Range("H2").Value = $latestHoDate
$de.Range("J2").Value = $newDeDeXlf
$de.Range("K2").Value = $deHandbackDate

# Row 3 mirrors row 2's handoff/handback file + datetimes (matches source data)
$de.Range("G3").Value = $newDeDeXlf
$de.Range("H3").Value = $latestHoDate
$de.Range("J3").Value = $newDeDeXlf
$de.Range("K3").Value = $deHandbackDate

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/89a8524106103f1b2796df6de832e7068317d359/e2e/$newMd1", [Type]::Missing, [Type]::Missing, $newMd1)
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3d8b721ce7931602742498f93f3c801e4faeced6/e2e/$newMd1", [Type]::Missing, [Type]::Missing, $newMd1)
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/89a8524106103f1b2796df6de832e7068317d359/e2e/$newMd2", [Type]::Missing, [Type]::Missing, $newMd2)
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3d8b721ce7931602742498f93f3c801e4faeced6/e2e/$newMd2", [Type]::Missing, [Type]::Missing, $newMd2)

Write-Output "Done."
